$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data rows for account 005428871 (ROSANGELA, 15871.44)
# and account 003512801 (LAIS, 3115.52) from the "Export" sheet.
$rowRosangela = $ws.Columns.Item(1).Find("005428871").Row
$ws.Rows.Item($rowRosangela).Delete()

$rowLais = $ws.Columns.Item(1).Find("003512801").Row
$ws.Rows.Item($rowLais).Delete()
